# Scheduled runner update: refresh Leve profit calculations (currentAveragePrice*
# columns and derived Leve price/profit columns) across the server sheets with
# the latest market-board pricing data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 847.2778
$ws.Range("I19").Value = 699.8889
$ws.Range("J19").Value = 994.6667
$ws.Range("K19").Value = 699.8889
$ws.Range("L19").Value = 994.6667
$ws.Range("M19").Value = -524.8889
$ws.Range("N19").Value = -1344.6667

$ws.Range("H39").Value = 138.55556
$ws.Range("I39").Value = 87.833336
$ws.Range("K39").Value = 263.500008
$ws.Range("M39").Value = 32.49999200000002

$ws.Range("H43").Value = 737.6
$ws.Range("I43").Value = 650
$ws.Range("J43").Value = 775.1429000000001
$ws.Range("K43").Value = 650
$ws.Range("L43").Value = 775.1429000000001
$ws.Range("M43").Value = -581
$ws.Range("N43").Value = -913.1429000000001

$ws.Range("H62").Value = 5088.0557
$ws.Range("J62").Value = 6167.778
$ws.Range("L62").Value = 6167.778
$ws.Range("N62").Value = -7415.778

$ws.Range("H65").Value = 5088.0557
$ws.Range("J65").Value = 6167.778
$ws.Range("L65").Value = 30838.89
$ws.Range("N65").Value = -37078.89

$ws.Range("H132").Value = 368816.78
$ws.Range("I132").Value = 434848.75
$ws.Range("J132").Value = 60667.668
$ws.Range("K132").Value = 1304546.25
$ws.Range("L132").Value = 182003.004
$ws.Range("M132").Value = -1302016.25
$ws.Range("N132").Value = -187063.004

$ws.Range("H137").Value = 55557304
$ws.Range("I137").Value = 90910690
$ws.Range("J137").Value = 1986.1428
$ws.Range("K137").Value = 272732070
$ws.Range("L137").Value = 5958.428400000001
$ws.Range("M137").Value = -272729520
$ws.Range("N137").Value = -11058.4284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H52").Value = 40000
$ws.Range("J52").Value = 40000
$ws.Range("L52").Value = 40000
$ws.Range("N52").Value = -40636

$ws.Range("H132").Value = 1960.3959
$ws.Range("I132").Value = 1692.5264
$ws.Range("J132").Value = 2978.3
$ws.Range("K132").Value = 5077.5792
$ws.Range("L132").Value = 8934.900000000001
$ws.Range("M132").Value = -2547.5792
$ws.Range("N132").Value = -13994.9

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1982.1052
$ws.Range("I20").Value = 1869.1
$ws.Range("J20").Value = 2107.6667
$ws.Range("K20").Value = 1869.1
$ws.Range("L20").Value = 2107.6667
$ws.Range("M20").Value = -1622.1
$ws.Range("N20").Value = -2601.6667

$ws.Range("H80").Value = 285.3158
$ws.Range("I80").Value = 350
$ws.Range("J80").Value = 277.70587
$ws.Range("K80").Value = 350
$ws.Range("L80").Value = 277.70587
$ws.Range("M80").Value = 648
$ws.Range("N80").Value = -2273.70587

$ws.Range("H83").Value = 285.3158
$ws.Range("I83").Value = 350
$ws.Range("J83").Value = 277.70587
$ws.Range("K83").Value = 1750
$ws.Range("L83").Value = 1388.52935
$ws.Range("M83").Value = 3242
$ws.Range("N83").Value = -11372.52935

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 18200
$ws.Range("J29").Value = 18200
$ws.Range("L29").Value = 18200
$ws.Range("N29").Value = -18786

$ws.Range("H86").Value = 2712.5715
$ws.Range("I86").Value = 2663.3333
$ws.Range("J86").Value = 2749.5
$ws.Range("K86").Value = 2663.3333
$ws.Range("L86").Value = 2749.5
$ws.Range("M86").Value = -1540.3333
$ws.Range("N86").Value = -4995.5

$ws.Range("H89").Value = 2712.5715
$ws.Range("I89").Value = 2663.3333
$ws.Range("J89").Value = 2749.5
$ws.Range("K89").Value = 13316.6665
$ws.Range("L89").Value = 13747.5
$ws.Range("M89").Value = -7700.666499999999
$ws.Range("N89").Value = -24979.5

$ws.Range("H94").Value = 1073.3793
$ws.Range("I94").Value = 803.8333
$ws.Range("J94").Value = 1143.6957
$ws.Range("K94").Value = 803.8333
$ws.Range("L94").Value = 1143.6957
$ws.Range("M94").Value = -352.8333
$ws.Range("N94").Value = -2045.6957

$ws.Range("H132").Value = 2661.2173
$ws.Range("I132").Value = 2147.7144
$ws.Range("J132").Value = 3460
$ws.Range("K132").Value = 6443.1432
$ws.Range("L132").Value = 10380
$ws.Range("M132").Value = -3913.1432
$ws.Range("N132").Value = -15440

$ws.Range("H134").Value = 2496.8667
$ws.Range("I134").Value = 1251.5
$ws.Range("J134").Value = 3920.1428
$ws.Range("K134").Value = 3754.5
$ws.Range("L134").Value = 11760.4284
$ws.Range("M134").Value = -1219.5
$ws.Range("N134").Value = -16830.4284

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 32.5
$ws.Range("I12").Value = 6.3333335
$ws.Range("J12").Value = 84.833336
$ws.Range("K12").Value = 19.0000005
$ws.Range("L12").Value = 254.500008
$ws.Range("M12").Value = 153.9999995
$ws.Range("N12").Value = -600.500008

$ws.Range("H16").Value = 645.5
$ws.Range("I16").Value = 645.5
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1936.5
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1763.5
$ws.Range("N16").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1497.8422
$ws.Range("I102").Value = 1218.25
$ws.Range("K102").Value = 1218.25
$ws.Range("M102").Value = 403.75

$ws.Range("H122").Value = 1236535.4
$ws.Range("I122").Value = 1236535.4
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3709606.2
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3707156.2
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 2888.2856
$ws.Range("I132").Value = 1723.0769
$ws.Range("J132").Value = 4781.75
$ws.Range("K132").Value = 5169.2307
$ws.Range("L132").Value = 14345.25
$ws.Range("M132").Value = -2639.2307
$ws.Range("N132").Value = -19405.25

$ws.Range("H137").Value = 44000
$ws.Range("J137").Value = 44000
$ws.Range("L137").Value = 44000
$ws.Range("N137").Value = -54200

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 6207.737
$ws.Range("I136").Value = 2250
$ws.Range("J136").Value = 9086.091
$ws.Range("K136").Value = 6750
$ws.Range("L136").Value = 27258.273
$ws.Range("M136").Value = -4200
$ws.Range("N136").Value = -32358.273

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 64100.312
$ws.Range("I126").Value = 92209.55
$ws.Range("J126").Value = 2260
$ws.Range("K126").Value = 276628.65
$ws.Range("L126").Value = 6780
$ws.Range("M126").Value = -274158.65
$ws.Range("N126").Value = -11720

$ws.Range("H136").Value = 22291120
$ws.Range("I136").Value = 30395304
$ws.Range("J136").Value = 4610.5
$ws.Range("K136").Value = 22291120
$ws.Range("L136").Value = 13831.5
$ws.Range("M136").Value = -91183362
$ws.Range("N136").Value = -18931.5
